# "Summary" sheet of the Fiji MSME indicators workbook.
# The two rows that list
#   - "Enterprises (absolute #)"                  / "4200"
#   - "Enterprises density (per 1000 people)"      / "4.8"
# need to swap order: the "density" row now comes first (row 11), followed
# by the "absolute #" row (row 12). Column A holds the label, column D
# holds the value (stored as text, not a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 -> becomes the "Enterprises density (per 1000 people)" / "4.8" row
$ws.Range("A11").Value2 = "Enterprises density (per 1000 people)"
$ws.Range("D11").Formula = "'4.8"
$ws.Range("D11").Style = "Normal"

# Row 12 -> becomes the "Enterprises (absolute #)" / "4200" row
$ws.Range("A12").Value2 = "Enterprises (absolute #)"
$ws.Range("D12").Formula = "'4200"
$ws.Range("D12").Style = "Normal"
